# Update "想去人数" (interested-attendee count) values that changed upstream.
# Row 3  (id=81916, 环形宇宙动漫游戏嘉年华): 8154 -> 8155
# Row 10 (id=82924, 首届运动番only):        198  -> 199
# These figures live on both the "展览" sheet and the "全部类型" sheet,
# which mirror the same underlying data.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 8155
    $ws.Range("F10").Value = 199
}
